$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 so the existing rows 5-8 shift down to 6-9
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly data point
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44469
$ws.Range("D5").Style = $ws.Range("D6").Style
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 100114007
$ws.Range("G5").Value = "Jengibre"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 140
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 13500
$ws.Range("N5").Value = "$/caja 13 kilos"
$ws.Range("O5").Value = "Perú"
$ws.Range("P5").Value = 1038
$ws.Range("Q5").Value = 13
$ws.Range("R5").Value = "Hortaliza"
